$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 158218
$ws.Range("C4").Value = 149274
$ws.Range("C7").Value = 5.65
$ws.Range("C8").Value = 64.04000000000001
